$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for species that were dropped entirely from the report.
# (Delete from the bottom up so earlier row numbers stay valid while deleting.)
$ws.Rows(53).Delete()   # Wood NA
$ws.Rows(48).Delete()   # Stones NA
$ws.Rows(46).Delete()   # Shells NA
$ws.Rows(43).Delete()   # Liocarcinus depurator
$ws.Rows(40).Delete()   # Eggs of Murex
$ws.Rows(37).Delete()   # Biological discard

# Zero out the W(kg) column for the remaining benthos rows (32-47), which now
# report weight pending a later re-measurement - except Dardanus calidus,
# which gets an updated weight of 0.014.
for ($r = 32; $r -le 47; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}
$ws.Cells.Item(38, 7).Value = 0.014
